# Update the example data
#  - Sheet1: rename the "Notes N" values in column C to "Note N"
#  - Sheet2: regenerate the age / sex / edu example data, and drop the
#    bold/centered header style so the header row uses the default style
#  - Leave Sheet1 as the active sheet with F6 selected (matches the
#    post-edit view state)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1: "Notes N" -> "Note N" (column C, rows 2..51)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$notes = @(
    "Note 1","Note 2","Note 3","Note 4","Note 5","Note 6","Note 7","Note 8","Note 9","Note 10",
    "Note 11","Note 12","Note 13","Note 14","Note 15","Note 16","Note 17","Note 18","Note 19","Note 20",
    "Note 21","Note 22","Note 23","Note 24","Note 25","Note 26","Note 27","Note 28","Note 29","Note 30",
    "Note 31","Note 32","Note 33","Note 34","Note 35","Note 36","Note 37","Note 38","Note 39","Note 40",
    "Note 41","Note 42","Note 43","Note 44","Note 45","Note 46","Note 47","Note 48","Note 49","Note 50"
)

for ($i = 0; $i -lt 50; $i++) {
    $ws1.Cells.Item($i + 2, 3).Value = $notes[$i]
}

# ---------------------------------------------------------------------------
# Sheet2: refresh the example age / sex / edu data
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Drop the bold / centered header style -> back to the workbook default
$ws2.Range("A1:D1").ClearFormats()

$ages = @(
    59,71,59,66,56,70,53,60,51,78,
    63,55,76,67,51,76,51,77,65,62,
    61,77,73,73,57,71,61,68,60,65,
    75,66,77,76,63,62,63,57,50,77,
    72,65,50,72,59,51,54,77,61,55
)

$sexes = @(
    "Male","Male","Female","Female","Female","Female","Male","Female","Male","Female",
    "Male","Male","Female","Female","Female","Female","Male","Male","Male","Male",
    "Male","Female","Male","Female","Female","Female","Female","Male","Female","Male",
    "Male","Male","Female","Female","Male","Male","Female","Female","Female","Male",
    "Male","Male","Female","Male","Male","Female","Female","Female","Male","Male"
)

$edus = @(
    5,7,10,12,8,8,10,7,9,6,
    10,5,5,7,7,14,10,8,5,10,
    5,8,9,9,7,6,11,5,5,6,
    5,6,10,7,6,7,12,9,13,8,
    5,6,8,6,6,10,5,10,10,8
)

for ($i = 0; $i -lt 50; $i++) {
    $r = $i + 2
    $ws2.Cells.Item($r, 2).Value = $ages[$i]
    $ws2.Cells.Item($r, 3).Value = $sexes[$i]
    $ws2.Cells.Item($r, 4).Value = $edus[$i]
}

# ---------------------------------------------------------------------------
# View state: Sheet1 becomes the active sheet, F6 selected
# ---------------------------------------------------------------------------
[void]$ws1.Activate()
[void]$ws1.Range("F6").Select()
